$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = 925
$ws.Range("E4").Value = 2
$ws.Range("C5").Value = 867
$ws.Range("E5").Value = 10
$ws.Range("F5").Value = 3.21
$ws.Range("C6").Value = 921
$ws.Range("F6").Value = 2.976
$ws.Range("C11").Value = 228
$ws.Range("C12").Value = 79
$ws.Range("E12").Value = 3
$ws.Range("C28").Value = 915
$ws.Range("C29").Value = 932
$ws.Range("E29").Value = 1
$ws.Range("C30").Value = 899
$ws.Range("E30").Value = 2
$ws.Range("F30").Value = 3.227
$ws.Range("F31").Value = 3.001
$ws.Range("C32").Value = 693
$ws.Range("C33").Value = 548
$ws.Range("E33").Value = 5
$ws.Range("F33").Value = 2.268
$ws.Range("C34").Value = 388
$ws.Range("E34").Value = 5
$ws.Range("F34").Value = 1.542
$ws.Range("C35").Value = 221
$ws.Range("E35").Value = 6
$ws.Range("F35").Value = 0.846
$ws.Range("C36").Value = 78
$ws.Range("E36").Value = 5
$ws.Range("F36").Value = 0.323
$ws.Range("E37").Value = 4
$ws.Range("F37").Value = 0.063
$ws.Range("E38").Value = 2
$ws.Range("E39").Value = 1
$ws.Range("E40").Value = 1
$ws.Range("E41").Value = 1
$ws.Range("C82").Value = 392
$ws.Range("C83").Value = 226
$ws.Range("E83").Value = 1
$ws.Range("F83").Value = 0.826
$ws.Range("E84").Value = 2
$ws.Range("F84").Value = 0.286
$ws.Range("C94").Value = 80
$ws.Range("E94").Value = 3
$ws.Range("C95").Value = 218
$ws.Range("E95").Value = 6
$ws.Range("F95").Value = 0.101
$ws.Range("C96").Value = 368
$ws.Range("E96").Value = 8
$ws.Range("F96").Value = 0.468
$ws.Range("C97").Value = 479
$ws.Range("E97").Value = 15
$ws.Range("F97").Value = 1.071
$ws.Range("C98").Value = 558
$ws.Range("E98").Value = 21
$ws.Range("F98").Value = 1.64
$ws.Range("C99").Value = 635
$ws.Range("E99").Value = 22
$ws.Range("F99").Value = 1.906
$ws.Range("C100").Value = 688
$ws.Range("E100").Value = 22
$ws.Range("F100").Value = 2.214
$ws.Range("C101").Value = 733
$ws.Range("E101").Value = 19
$ws.Range("F101").Value = 2.4
$ws.Range("C102").Value = 759
$ws.Range("E102").Value = 14
$ws.Range("F102").Value = 2.485
$ws.Range("C103").Value = 738
$ws.Range("E103").Value = 9
$ws.Range("F103").Value = 2.647
$ws.Range("C104").Value = 653
$ws.Range("E104").Value = 7
$ws.Range("F104").Value = 2.464
$ws.Range("C105").Value = 528
$ws.Range("E105").Value = 5
$ws.Range("C106").Value = 374
$ws.Range("E106").Value = 4
$ws.Range("F106").Value = 1.418
$ws.Range("C107").Value = 207
$ws.Range("F107").Value = 0.84
$ws.Range("E108").Value = 7
$ws.Range("F108").Value = 0.311
$ws.Range("E109").Value = 7
$ws.Range("F109").Value = 0.056
$ws.Range("E110").Value = 7
$ws.Range("E111").Value = 7
$ws.Range("E112").Value = 6
$ws.Range("E113").Value = 3
$ws.Range("E117").Value = 0
$ws.Range("C118").Value = 83
$ws.Range("E118").Value = 0
$ws.Range("C119").Value = 230
$ws.Range("E119").Value = 0
$ws.Range("F119").Value = 0.109
$ws.Range("C120").Value = 396
$ws.Range("E120").Value = 0
$ws.Range("F120").Value = 0.442
$ws.Range("C121").Value = 557
$ws.Range("E121").Value = 1
$ws.Range("F121").Value = 1.134
$ws.Range("C122").Value = 700
$ws.Range("E122").Value = 1
$ws.Range("F122").Value = 1.836
$ws.Range("C123").Value = 813
$ws.Range("F123").Value = 2.394
$ws.Range("C124").Value = 881
$ws.Range("C125").Value = 901
$ws.Range("E125").Value = 1
$ws.Range("F125").Value = 3.085
$ws.Range("C126").Value = 877
$ws.Range("E126").Value = 1
$ws.Range("F126").Value = 2.945
$ws.Range("C127").Value = 806
$ws.Range("E127").Value = 1
$ws.Range("F127").Value = 2.859
$ws.Range("C128").Value = 697
$ws.Range("F128").Value = 2.586
$ws.Range("F129").Value = 2.263
$ws.Range("C148").Value = 878
$ws.Range("C149").Value = 890
$ws.Range("E149").Value = 2
$ws.Range("C150").Value = 838
$ws.Range("E150").Value = 4
$ws.Range("F150").Value = 2.914
$ws.Range("C151").Value = 751
$ws.Range("E151").Value = 7
$ws.Range("F151").Value = 2.722
$ws.Range("C152").Value = 648
$ws.Range("E152").Value = 7
$ws.Range("F152").Value = 2.501
$ws.Range("C153").Value = 517
$ws.Range("E153").Value = 6
$ws.Range("F153").Value = 2.025
$ws.Range("C154").Value = 368
$ws.Range("E154").Value = 4
$ws.Range("F154").Value = 1.311
$ws.Range("C155").Value = 213
$ws.Range("E155").Value = 2
$ws.Range("F155").Value = 0.84
$ws.Range("C156").Value = 73
$ws.Range("E156").Value = 1
$ws.Range("F156").Value = 0.262
$ws.Range("F157").Value = 0.056
$ws.Range("E159").Value = 2
$ws.Range("E160").Value = 6
$ws.Range("E161").Value = 11
$ws.Range("E162").Value = 11
$ws.Range("E163").Value = 12
$ws.Range("E164").Value = 16
$ws.Range("C165").Value = 3
$ws.Range("E165").Value = 18
$ws.Range("C166").Value = 68
$ws.Range("E166").Value = 16
$ws.Range("C167").Value = 189
$ws.Range("E167").Value = 16
$ws.Range("F167").Value = 0.152
$ws.Range("C168").Value = 322
$ws.Range("E168").Value = 18
$ws.Range("F168").Value = 0.461
$ws.Range("C169").Value = 481
$ws.Range("E169").Value = 14
$ws.Range("F169").Value = 0.917
$ws.Range("C170").Value = 636
$ws.Range("E170").Value = 9
$ws.Range("F170").Value = 1.617

Write-Host "Applied 166 cell updates"
